$p = $ppt.ActivePresentation
$s2 = $p.Slides.Item(2)

# --- Change 2 (done first, before slide 2's runs get merged): add a new
# slide 3, a duplicate of slide 2 (still with "2", ". ", "commit" as 3
# separate runs), then retarget the first run's text to "3" ---
$new = $s2.Duplicate().Item(1)
$shp3 = $new.Shapes.Item(1)
$tr3 = $shp3.TextFrame.TextRange
$tr3.Characters(1, 1).Text = "3"

# --- Change 1: slide 2 - merge the ". " and "commit" runs into ". commit" ---
$shp2 = $s2.Shapes.Item(1)
$tr2 = $shp2.TextFrame.TextRange
$tr2.Characters(2, 8).Text = ". commit"
